$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.6164125
$ws.Range("H2").Value = 3.232825
$ws.Range("M2").Value = 0.690242
$ws.Range("N2").Value = 1.380484
$ws.Range("O2").Value = 0.07603312339759918
$ws.Range("P2").Value = 0.05971740901392503
$ws.Range("Q2").Value = 1.115715796825
$ws.Range("R2").Value = 4.4628631873
$ws.Range("S2").Value = 0.07603312339759918
$ws.Range("T2").Value = 0.05971740901392503

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.6164125
$ws.Range("H3").Value = 3.232825
$ws.Range("O3").Value = 0.1680093924898635
$ws.Range("P3").Value = 0.1979351595954898
$ws.Range("Q3").Value = 2.465382518033334
$ws.Range("R3").Value = 14.7922951082
$ws.Range("S3").Value = 0.1680093924898635
$ws.Range("T3").Value = 0.1979351595954898

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.6164125
$ws.Range("H4").Value = 3.232825
$ws.Range("M4").Value = 1.488848666666667
$ws.Range("N4").Value = 4.466546
$ws.Range("O4").Value = 0.1640030806485518
$ws.Range("P4").Value = 0.1932152450600737
$ws.Range("Q4").Value = 2.406593595408333
$ws.Range("R4").Value = 14.43956157245
$ws.Range("S4").Value = 0.1640030806485518
$ws.Range("T4").Value = 0.1932152450600737

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.6164125
$ws.Range("H5").Value = 3.232825
$ws.Range("M5").Value = 3.427339
$ws.Range("N5").Value = 6.854678
$ws.Range("O5").Value = 0.3775361237253082
$ws.Range("P5").Value = 0.296521806688635
$ws.Range("Q5").Value = 5.5399936013375
$ws.Range("R5").Value = 22.15997440535
$ws.Range("S5").Value = 0.3775361237253082
$ws.Range("T5").Value = 0.296521806688635

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.6164125
$ws.Range("H6").Value = 3.232825
$ws.Range("M6").Value = 1.086196666666667
$ws.Range("N6").Value = 3.25859
$ws.Range("O6").Value = 0.1196492319950504
$ws.Range("P6").Value = 0.1409611062777156
$ws.Range("Q6").Value = 1.755741869458333
$ws.Range("R6").Value = 10.53445121675
$ws.Range("S6").Value = 0.1196492319950504
$ws.Range("T6").Value = 0.1409611062777156

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.6164125
$ws.Range("H7").Value = 3.232825
$ws.Range("M7").Value = 0.8603300000000002
$ws.Range("N7").Value = 2.58099
$ws.Range("O7").Value = 0.09476904774362691
$ws.Range("P7").Value = 0.111649273364161
$ws.Range("Q7").Value = 1.390648166125
$ws.Range("R7").Value = 8.343888996750001
$ws.Range("S7").Value = 0.09476904774362691
$ws.Range("T7").Value = 0.111649273364161
